# Clear the "No. of Sites/bldg ..." breakdown columns (AB:AK) and the
# DIFFERENCE column (AM) for all data rows, keeping PREVIOUS ACCOMPLISHMENT
# (AL) intact, per the "most updated status accomplishment files as of may"
# commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 19; $row++) {
    $ws.Range("AB$row`:AK$row").ClearContents()
    $ws.Range("AM$row").ClearContents()
}
